# Updates cryptos list data (prices / 1h volume change %) to match the
# latest scrape, as described in the commit message:
# "Updated cryptos list on Thu Feb  8 07:47:21 UTC 2024 with GitHub Actions"
#
# All cells in this sheet hold plain text (coin names, URLs, price strings
# that sometimes use dotted thousands separators, and padded percentage
# strings). Assigning directly via Range.Value lets Excel's usual
# type-inference re-interpret simple decimal-looking strings (e.g. "1.00"
# or "6.94") as numbers, which would silently change the stored cell type.
# To keep every touched cell as literal text (matching the source data),
# we briefly force Text number-formatting for the assignment and then
# restore the cell to the default "Normal" style so no formatting residue
# is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "44.515.74"
Set-TextValue "E2" "  +3.68%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.421.54"
Set-TextValue "E3" "  +2.52%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "314.50"
Set-TextValue "E5" "  +4.08%  "

# Row 6 - Solana
Set-TextValue "D6" "100.94"
Set-TextValue "E6" "  +5.85%  "

# Row 7 - XRP
Set-TextValue "E7" "  +2.34%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.518"
Set-TextValue "E9" "  +6.52%  "

# Row 10 - Avalanche
Set-TextValue "D10" "35.40"
Set-TextValue "E10" "  +4.12%  "

# Row 11 and Row 12 swap places: Chainlink <-> Dogecoin
Set-TextValue "B11" "Dogecoin"
Set-TextValue "C11" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D11" "0.0799"
Set-TextValue "E11" "  +1.82%  "

Set-TextValue "B12" "Chainlink"
Set-TextValue "C12" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D12" "19.17"
Set-TextValue "E12" "  +4.56%  "

# Row 13 - TRON
Set-TextValue "E13" "  -2.07%  "

# Row 14 - Polkadot
Set-TextValue "D14" "6.94"
Set-TextValue "E14" "  +3.56%  "

# Row 15 - Wrapped liquid staked Ether 2.0
Set-TextValue "D15" "2.800.76"
Set-TextValue "E15" "  +2.62%  "

# Row 16 - Wrapped Ether
Set-TextValue "D16" "2.397.92"
Set-TextValue "E16" "  +0.84%  "

# Row 17 - Polygon
Set-TextValue "E17" "  +5.39%  "

# Row 18 - Wrapped BTC
Set-TextValue "D18" "44.342.92"
Set-TextValue "E18" "  +3.36%  "

# Row 19 - Internet Computer (DFINITY)
Set-TextValue "D19" "12.41"
Set-TextValue "E19" "  +4.32%  "

# Row 20 - Uniswap
Set-TextValue "D20" "6.40"
Set-TextValue "E20" "  +2.35%  "

# Row 21 - Shiba Inu
Set-TextValue "D21" "0.0₃0918"
Set-TextValue "E21" "  +3.64%  "

# Row 22 - Litecoin
Set-TextValue "E22" "  +1.02%  "

# Row 23 - Bitcoin Cash
Set-TextValue "D23" "241.74"
Set-TextValue "E23" "  +2.80%  "

# Row 24 - Immutable X
Set-TextValue "E24" "  +5.64%  "

# Row 25 - PancakeSwap
Set-TextValue "E25" "  +1.94%  "

# Row 26 - Dai
Set-TextValue "E26" "  -0.04%  "

# Row 27 - Ethereum Classic
Set-TextValue "D27" "25.18"
Set-TextValue "E27" "  +2.98%  "

# Row 28 - Toncoin
Set-TextValue "E28" "  -3.71%  "

# Row 29 - Cosmos
Set-TextValue "E29" "  +3.07%  "

# Row 30 - Injective Protocol
Set-TextValue "D30" "33.23"
Set-TextValue "E30" "  +4.16%  "

# Row 31 - OKB
Set-TextValue "E31" "  +0.82%  "

# Row 32 - Kaspa
Set-TextValue "E32" "  +17.50%  "

# Row 33 - Celestia
Set-TextValue "D33" "19.29"
Set-TextValue "E33" "  +10.27%  "

# Row 34 - Filecoin
Set-TextValue "E34" "  +3.39%  "

# Row 35 - Hedera
Set-TextValue "E35" "  +7.59%  "

# Row 36 - First Digital USD
Set-TextValue "E36" "  +0.19%  "

# Row 37 - ARBITRUM
Set-TextValue "E37" "  +2.54%  "

# Row 38 - Render Token
Set-TextValue "E38" "  +3.45%  "

# Row 39 - Lido DAO Token
Set-TextValue "E39" "  +1.94%  "

# Row 40 - WEMIX Token
Set-TextValue "E40" "  -1.93%  "

# Row 41 - Monero
Set-TextValue "D41" "121.35"
Set-TextValue "E41" "  -7.11%  "

# Row 42 - Stellar
Set-TextValue "E42" "  +1.48%  "

# Row 43 - EnergySwap
Set-TextValue "D43" "20.85"
Set-TextValue "E43" "  -1.86%  "

# Row 44 - VeChain
Set-TextValue "E44" "  +3.84%  "

# Row 45 - Maker
Set-TextValue "D45" "1.943.71"
Set-TextValue "E45" "  +0.73%  "

# Row 46 - ApeX Protocol
Set-TextValue "E46" "  +1.93%  "

# Row 47 - NEAR Protocol
Set-TextValue "D47" "2.93"
Set-TextValue "E47" "  +8.54%  "

# Row 48 - Frax Share
Set-TextValue "E48" "  +3.23%  "

# Row 49 - Stacks
Set-TextValue "D49" "1.66"
Set-TextValue "E49" "  +10.04%  "

# Row 50 - MultiversX
Set-TextValue "D50" "55.05"
Set-TextValue "E50" "  +7.11%  "

# Row 51 - Bitcoin SV
Set-TextValue "D51" "74.54"
Set-TextValue "E51" "  +4.30%  "
